$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values, re-pulled from source data
$values = @{
    2  = -1
    3  = 3
    4  = 1
    5  = -5
    6  = 1
    7  = -3
    9  = 5
    10 = -3
    11 = -1
    12 = -4
    13 = 2
    14 = -2
    15 = -1
    16 = 4
    17 = -3
    18 = 8
    19 = 3
    20 = -1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
